# Scheduled runner update: refresh currentAveragePrice / leve profit figures
# derived from the latest market-board snapshot (per-sheet leve tables).
$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1927.5714
$ws.Range("J19").Value = 455.83334
$ws.Range("L19").Value = 455.83334
$ws.Range("N19").Value = -805.83334
$ws.Range("H68").Value = 54647.332
$ws.Range("I68").Value = 60000
$ws.Range("K68").Value = 60000
$ws.Range("M68").Value = -59251
$ws.Range("H71").Value = 54647.332
$ws.Range("I71").Value = 60000
$ws.Range("K71").Value = 180000
$ws.Range("M71").Value = -176256
$ws.Range("H132").Value = 1421.3182
$ws.Range("I132").Value = 1414.05
$ws.Range("J132").Value = 1494
$ws.Range("K132").Value = 4242.15
$ws.Range("L132").Value = 4482
$ws.Range("M132").Value = -1712.15
$ws.Range("N132").Value = -9542
$ws.Range("H138").Value = 4888.3335
$ws.Range("I138").Value = 4000
$ws.Range("J138").Value = 4999.375
$ws.Range("K138").Value = 12000
$ws.Range("L138").Value = 14998.125
$ws.Range("M138").Value = -6860
$ws.Range("N138").Value = -25278.125

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 9193.474
$ws.Range("I61").Value = 3471.5386
$ws.Range("J61").Value = 21591
$ws.Range("K61").Value = 3471.5386
$ws.Range("L61").Value = 21591
$ws.Range("M61").Value = -3259.5386
$ws.Range("N61").Value = -22015
$ws.Range("H75").Value = 50083.332
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 50083.332
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 50083.332
$ws.Range("M75").ClearContents()
$ws.Range("N75").Value = -51831.332
$ws.Range("H78").Value = 50083.332
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 50083.332
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 150249.996
$ws.Range("M78").ClearContents()
$ws.Range("N78").Value = -158985.996
$ws.Range("H132").Value = 7153.476
$ws.Range("I132").Value = 5575
$ws.Range("J132").Value = 22149
$ws.Range("K132").Value = 16725
$ws.Range("L132").Value = 66447
$ws.Range("M132").Value = -14195
$ws.Range("N132").Value = -71507
$ws.Range("H136").Value = 9193.474
$ws.Range("I136").Value = 3471.5386
$ws.Range("J136").Value = 21591
$ws.Range("K136").Value = 10414.6158
$ws.Range("L136").Value = 64773
$ws.Range("M136").Value = -7864.6158
$ws.Range("N136").Value = -69873

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1994.3125
$ws.Range("I107").Value = 1762.2307
$ws.Range("K107").Value = 1762.2307
$ws.Range("M107").Value = 157.7692999999999

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 36257.55
$ws.Range("I99").Value = 6496.5
$ws.Range("J99").Value = 41837.75
$ws.Range("K99").Value = 6496.5
$ws.Range("L99").Value = 41837.75
$ws.Range("M99").Value = -4998.5
$ws.Range("N99").Value = -44833.75
$ws.Range("H122").Value = 2116.2856
$ws.Range("I122").Value = 2326.5
$ws.Range("K122").Value = 6979.5
$ws.Range("M122").Value = -4529.5
$ws.Range("H126").Value = 36257.55
$ws.Range("I126").Value = 6496.5
$ws.Range("J126").Value = 41837.75
$ws.Range("K126").Value = 19489.5
$ws.Range("L126").Value = 125513.25
$ws.Range("M126").Value = -17019.5
$ws.Range("N126").Value = -130453.25
$ws.Range("H132").Value = 20296.986
$ws.Range("I132").Value = 14713.387
$ws.Range("K132").Value = 44140.161
$ws.Range("M132").Value = -41610.161
$ws.Range("H134").Value = 3591.2983
$ws.Range("I134").Value = 2257.2173
$ws.Range("J134").Value = 9170.182000000001
$ws.Range("K134").Value = 6771.651899999999
$ws.Range("L134").Value = 27510.546
$ws.Range("M134").Value = -4236.651899999999
$ws.Range("N134").Value = -32580.546

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 47.875
$ws.Range("J2").Value = 42.333332
$ws.Range("L2").Value = 253.999992
$ws.Range("N2").Value = -479.999992
$ws.Range("H107").Value = 1274.3043
$ws.Range("J107").Value = 1216
$ws.Range("L107").Value = 3648
$ws.Range("N107").Value = -7488

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2342.3333
$ws.Range("I102").Value = 2216.8948
$ws.Range("J102").Value = 2640.25
$ws.Range("K102").Value = 2216.8948
$ws.Range("L102").Value = 2640.25
$ws.Range("M102").Value = -594.8948
$ws.Range("N102").Value = -5884.25
$ws.Range("H126").Value = 3140.0386
$ws.Range("I126").Value = 3211.7368
$ws.Range("J126").Value = 2945.4285
$ws.Range("K126").Value = 9635.2104
$ws.Range("L126").Value = 8836.2855
$ws.Range("M126").Value = -7165.2104
$ws.Range("N126").Value = -13776.2855
$ws.Range("H132").Value = 16389.47
$ws.Range("I132").Value = 20459
$ws.Range("J132").Value = 6622.6
$ws.Range("K132").Value = 61377
$ws.Range("L132").Value = 19867.8
$ws.Range("M132").Value = -58847
$ws.Range("N132").Value = -24927.8

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 478.5
$ws.Range("I22").Value = 495.1111
$ws.Range("J22").Value = 428.66666
$ws.Range("K22").Value = 495.1111
$ws.Range("L22").Value = 428.66666
$ws.Range("M22").Value = -200.1111
$ws.Range("N22").Value = -1018.66666
$ws.Range("H27").Value = 478.5
$ws.Range("I27").Value = 495.1111
$ws.Range("J27").Value = 428.66666
$ws.Range("K27").Value = 495.1111
$ws.Range("L27").Value = 428.66666
$ws.Range("M27").Value = -388.1111
$ws.Range("N27").Value = -642.66666
$ws.Range("H68").Value = 2485.075
$ws.Range("I68").Value = 2138.4707
$ws.Range("J68").Value = 4449.1665
$ws.Range("K68").Value = 2138.4707
$ws.Range("L68").Value = 4449.1665
$ws.Range("M68").Value = -1389.4707
$ws.Range("N68").Value = -5947.1665
$ws.Range("H71").Value = 2485.075
$ws.Range("I71").Value = 2138.4707
$ws.Range("J71").Value = 4449.1665
$ws.Range("K71").Value = 10692.3535
$ws.Range("L71").Value = 22245.8325
$ws.Range("M71").Value = -6948.353499999999
$ws.Range("N71").Value = -29733.8325
$ws.Range("H88").Value = 54900
$ws.Range("J88").Value = 54900
$ws.Range("L88").Value = 54900
$ws.Range("N88").Value = -55756
$ws.Range("H91").Value = 54900
$ws.Range("J91").Value = 54900
$ws.Range("L91").Value = 54900
$ws.Range("N91").Value = -57864
$ws.Range("H122").Value = 3565.4443
$ws.Range("I122").Value = 3473.625
$ws.Range("K122").Value = 10420.875
$ws.Range("M122").Value = -7970.875
$ws.Range("H132").Value = 6495.2705
$ws.Range("I132").Value = 5871.8
$ws.Range("J132").Value = 9167.286
$ws.Range("K132").Value = 17615.4
$ws.Range("L132").Value = 27501.858
$ws.Range("M132").Value = -15085.4
$ws.Range("N132").Value = -32561.858

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 4879.3687
$ws.Range("I122").Value = 3479.5
$ws.Range("J122").Value = 7279.143
$ws.Range("K122").Value = 10438.5
$ws.Range("L122").Value = 21837.429
$ws.Range("M122").Value = -7988.5
$ws.Range("N122").Value = -26737.429
$ws.Range("H132").Value = 13308.737
$ws.Range("I132").Value = 6986.4756
$ws.Range("J132").Value = 33606.527
$ws.Range("K132").Value = 20959.4268
$ws.Range("L132").Value = 100819.581
$ws.Range("M132").Value = -18429.4268
$ws.Range("N132").Value = -105879.581
